# Applies the "Updated symbol list" data refresh to the crypto price sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $addr, $val) {
    # Force the cell to be stored as text so numeric-looking strings
    # (prices, percentages, hour numbers) are not coerced into numbers.
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

$rowUpdates = @(
    @{ "Row"=2; "D"="298.64"; "E"="0.64%"; "G"="12" }
    @{ "Row"=3; "D"="31.31"; "E"="-0.19%"; "G"="12" }
    @{ "Row"=4; "D"="5.098"; "E"="-0.27%"; "G"="12" }
    @{ "Row"=5; "D"="0.08020"; "E"="9.38%"; "G"="12" }
    @{ "Row"=6; "D"="2.603"; "E"="47.65%"; "G"="12" }
    @{ "Row"=7; "D"="7.811"; "E"="1.17%"; "G"="12" }
    @{ "Row"=8; "D"="3.822"; "E"="2.64%"; "G"="12" }
    @{ "Row"=9; "D"="0.9178"; "E"="-0.45%"; "G"="12" }
    @{ "Row"=10; "E"="3.67%"; "G"="12" }
    @{ "Row"=11; "D"="0.07359"; "E"="4.15%"; "G"="12" }
    @{ "Row"=12; "D"="0.08385"; "E"="4.84%"; "G"="12" }
    @{ "Row"=13; "D"="0.03024"; "E"="1.15%"; "G"="12" }
    @{ "Row"=14; "D"="0.09962"; "E"="0.65%"; "G"="12" }
    @{ "Row"=15; "D"="0.001507"; "E"="0.96%"; "G"="12" }
    @{ "Row"=16; "D"="0.005933"; "E"="-5.13%"; "G"="12" }
    @{ "Row"=17; "D"="3.503"; "E"="1.53%"; "G"="12" }
    @{ "Row"=18; "E"="1.33%"; "G"="12" }
    @{ "Row"=19; "E"="0.41%"; "G"="12" }
    @{ "Row"=20; "E"="0.41%"; "G"="12" }
    @{ "Row"=21; "D"="4.585"; "E"="0.75%"; "G"="12" }
    @{ "Row"=22; "D"="0.1599"; "E"="3.26%"; "G"="12" }
    @{ "Row"=23; "D"="0.04606"; "E"="-0.87%"; "G"="12" }
    @{ "Row"=24; "D"="0.001240"; "E"="1.38%"; "G"="12" }
    @{ "Row"=25; "D"="0.004453"; "E"="-6.30%"; "G"="12" }
    @{ "Row"=26; "D"="0.0001190"; "E"="-8.35%"; "G"="12" }
    @{ "Row"=27; "D"="0.0003427"; "E"="83.03%"; "G"="12" }
    @{ "Row"=28; "G"="12" }
    @{ "Row"=29; "G"="12" }
    @{ "Row"=30; "G"="12" }
    @{ "Row"=31; "G"="12" }
    @{ "Row"=32; "G"="12" }
    @{ "Row"=33; "G"="12" }
    @{ "Row"=34; "G"="12" }
    @{ "Row"=35; "G"="12" }
    @{ "Row"=36; "G"="12" }
    @{ "Row"=37; "G"="12" }
    @{ "Row"=38; "G"="12" }
    @{ "Row"=39; "D"="0.01831"; "E"="7.37%"; "G"="12" }
    @{ "Row"=40; "D"="0.04544"; "E"="2.69%"; "G"="12" }
    @{ "Row"=41; "D"="0.007056"; "E"="-1.92%"; "G"="12" }
    @{ "Row"=42; "E"="1.05%"; "G"="12" }
    @{ "Row"=43; "D"="0.002240"; "E"="4.86%"; "G"="12" }
    @{ "Row"=44; "D"="0.009703"; "E"="-11.34%"; "G"="12" }
    @{ "Row"=45; "D"="0.00006531"; "E"="7.68%"; "G"="12" }
    @{ "Row"=46; "E"="-0.03%"; "G"="12" }
    @{ "Row"=47; "B"="BOLO"; "C"="https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"; "D"="0.8206"; "E"="-57.27%"; "G"="12" }
    @{ "Row"=48; "B"="CoinbaseStockToken"; "C"="https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"; "D"="0.006196"; "E"="-39.33%"; "G"="12" }
    @{ "Row"=49; "E"="-0.03%"; "G"="12" }
    @{ "Row"=50; "E"="0.04%"; "G"="12" }
    @{ "Row"=51; "G"="12" }
)

foreach ($u in $rowUpdates) {
    $r = $u["Row"]
    foreach ($col in "B","C") {
        if ($u.ContainsKey($col)) {
            $ws.Range("$col$r").Value = $u[$col]
        }
    }
    foreach ($col in "D","E","G") {
        if ($u.ContainsKey($col)) {
            Set-TextCell $ws "$col$r" $u[$col]
        }
    }
}
